# Edit script: add COOPERMIL cooperative data to master_resultados workbook
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # comparativo_master
$ws2 = $wb.Worksheets.Item(2)   # niveis_master
$ws3 = $wb.Worksheets.Item(3)   # financeiro_master
$ws4 = $wb.Worksheets.Item(4)   # questionario

# ---------------------------------------------------------------------------
# 1) comparativo_master: append 17 detail rows (170-186) for COOPERMIL
# ---------------------------------------------------------------------------
$newRows = @(
    "COOPERMIL|ADRIANO RAFAEL DILLY|60|Avançado|85|Avançado|25|41.666666666666671",
    "COOPERMIL|ANDRE LUCIANO RIECKE|40|Intermediário|69|Avançado|29|72.5",
    "COOPERMIL|ANDREA REGINA BRINCKER|52|Intermediário|75|Avançado|23|44.230769230769234",
    "COOPERMIL|ARMANDO PETRY|37|Intermediário|75|Avançado|38|102.70270270270269",
    "COOPERMIL|DIRCEU GEREMIA|30|Intermediário|60|Avançado|30|100",
    "COOPERMIL|EDUARDO ANDRE ULLMANN|29|Básico|58|Intermediário|29|100",
    "COOPERMIL|EDUARDO FACCHINELLO|21|Básico|55|Intermediário|34|161.9047619047619",
    "COOPERMIL|ELIZEU MAZZARRO|25|Básico|62|Avançado|37|148",
    "COOPERMIL|FELIPE GABRIEL GAVIRAGHI|34|Intermediário|62|Avançado|28|82.35294117647058",
    "COOPERMIL|GIOVANA PILECCO HERMANN|42|Intermediário|79|Avançado|37|88.095238095238088",
    "COOPERMIL|GUILHERME HENRIQUE WAGNER|27|Básico|81|Avançado|54|200",
    "COOPERMIL|IVAN ROBERTO HAAS|41|Intermediário|76|Avançado|35|85.365853658536579",
    "COOPERMIL|JAIRO LEANDRO MULLER|13|Básico|55|Intermediário|42|323.07692307692309",
    "COOPERMIL|JUAREZ ANDRE BECK|27|Básico|50|Intermediário|23|85.18518518518519",
    "COOPERMIL|MAURICIO ANDRE HORN|20|Básico|56|Intermediário|36|180",
    "COOPERMIL|OSVINO FRISKE|50|Intermediário|77|Avançado|27|54",
    "COOPERMIL|RICARDO ANTÔNIO PIZZONI|21|Básico|57|Intermediário|36|171.42857142857139"
)

$startRow = 170
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $parts = $newRows[$i].Split("|")
    $r = $startRow + $i
    $ws1.Cells.Item($r, 1).Value = $parts[0]
    $ws1.Cells.Item($r, 2).Value = $parts[1]
    $ws1.Cells.Item($r, 3).Value = $parts[2]
    $ws1.Cells.Item($r, 4).Value = $parts[3]
    $ws1.Cells.Item($r, 5).Value = $parts[4]
    $ws1.Cells.Item($r, 6).Value = $parts[5]
    $ws1.Cells.Item($r, 7).Value = $parts[6]
    $ws1.Cells.Item($r, 8).Value = $parts[7]
}
$ws1.Range("B170:H186").VerticalAlignment = -4108

$ws1.Columns.Item(3).ColumnWidth = 15.6640625
$ws1.Columns.Item(4).ColumnWidth = 13
$ws1.Columns.Item(5).ColumnWidth = 14.6640625
$ws1.Columns.Item(8).ColumnWidth = 13.44140625

# ---------------------------------------------------------------------------
# 2) niveis_master: insert 3 summary rows (32-34) for COOPERMIL by level
# ---------------------------------------------------------------------------
$ws2.Range("A32:A34").EntireRow.Insert()

$ws2.Range("A32").Value = "COOPERMIL"
$ws2.Range("B32").Value = "Básico"
$ws2.Range("C32").Value = 8
$ws2.Range("D32").Value = 0

$ws2.Range("A33").Value = "COOPERMIL"
$ws2.Range("B33").Value = "Intermediário"
$ws2.Range("C33").Value = 8
$ws2.Range("D33").Value = 6

$ws2.Range("A34").Value = "COOPERMIL"
$ws2.Range("B34").Value = "Avançado"
$ws2.Range("C34").Value = 1
$ws2.Range("D34").Value = 11

$ws2.Range("C32:D34").VerticalAlignment = -4108

$ws2.Columns.Item(1).ColumnWidth = 13.21875
$ws2.Columns.Item(2).ColumnWidth = 13
$ws2.Columns.Item(4).ColumnWidth = 8.77734375

# ---------------------------------------------------------------------------
# 3) financeiro_master: insert 1 summary row (12) for COOPERMIL
# ---------------------------------------------------------------------------
$ws3.Range("A12").EntireRow.Insert()

$ws3.Range("A12").Value = "COOPERMIL"
$ws3.Range("B12").Value = "Gestão Financeira"
$ws3.Range("C12").Value = 94
$ws3.Range("D12").Value = 200
$ws3.Range("E12").Value = 106
$ws3.Range("F12").Value = 112.7659574468085

$ws3.Columns.Item(1).ColumnWidth = 13.21875
$ws3.Columns.Item(2).ColumnWidth = 16.33203125
$ws3.Columns.Item(3).ColumnWidth = 18.21875
$ws3.Columns.Item(4).ColumnWidth = 17.21875
$ws3.Columns.Item(5).ColumnWidth = 16.88671875
$ws3.Columns.Item(6).ColumnWidth = 13.44140625

# ---------------------------------------------------------------------------
# 4) Selections / active sheet, matching the saved workbook view state
# ---------------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("G22").Select()

$ws3.Activate()
$ws3.Range("C12:F12").Select()

$ws4.Activate()
$ws4.Range("I5").Select()

$ws1.Activate()
$ws1.Range("B170:H186").Select()
